$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$win = $excel.ActiveWindow

# keep freeze as-is (already frozen 1,1 at B2)
$win.ScrollRow = 33
$win.ScrollColumn = 1
$ws.Range("I51").Select()
Write-Output $win.VisibleRange.Address()
